$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cw0")
Write-Output $ws.Name
